$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I: convert the "$x,xxx" text earnings into real numbers with a
# --- number format, mirroring what Excel does when you paste/convert a
# --- currency-text column into numeric currency values.

# Standard "#,##0" rows
$ws.Range("I2").Value  = 15000
$ws.Range("I3").Value  = 8000
$ws.Range("I4").Value  = 6000
$ws.Range("I5").Value  = 5000
$ws.Range("I6").Value  = 4000
$ws.Range("I7").Value  = 3000
$ws.Range("I9").Value  = 2400
$ws.Range("I10").Value = 2300
$ws.Range("I11").Value = 2200
$ws.Range("I12").Value = 2100
$ws.Range("I13").Value = 2000
$ws.Range("I14").Value = 1950
$ws.Range("I15").Value = 1900
$ws.Range("I16").Value = 1850
$ws.Range("I17").Value = 1800
$ws.Range("I18").Value = 1750
$ws.Range("I19").Value = 1700
$ws.Range("I20").Value = 1650
$ws.Range("I21").Value = 1600
$ws.Range("I22").Value = 1550
$ws.Range("I23").Value = 1500
$ws.Range("I24").Value = 1450
$ws.Range("I25").Value = 1400
$ws.Range("I26").Value = 1225
$ws.Range("I27").Value = 1225
$ws.Range("I28").Value = 1225
$ws.Range("I29").Value = 1225
$ws.Range("I30").Value = 1225
$ws.Range("I31").Value = 1225
$ws.Range("I32").Value = 1225
$ws.Range("I35").Value = 1150
$ws.Range("I36").Value = 1150
$ws.Range("I37").Value = 1150
$ws.Range("I38").Value = 1150
$ws.Range("I39").Value = 1150
$ws.Range("I40").Value = 1050
$ws.Range("I43").Value = 1050
$ws.Range("I44").Value = 1050
$ws.Range("I49").Value = 1050
$ws.Range("I53").Value = 1050
$ws.Range("I56").Value = 1050
$ws.Range("I59").Value = 1050

$ws.Range("I2:I7").NumberFormat   = "#,##0"
$ws.Range("I9:I32").NumberFormat  = "#,##0"
$ws.Range("I35:I40").NumberFormat = "#,##0"
$ws.Range("I43:I44").NumberFormat = "#,##0"
$ws.Range("I49").NumberFormat     = "#,##0"
$ws.Range("I53").NumberFormat     = "#,##0"
$ws.Range("I56").NumberFormat     = "#,##0"
$ws.Range("I59").NumberFormat     = "#,##0"

# Currency format row (kept visually "$2,500" style, with red negatives)
$ws.Range("I8").Value = 2500
$ws.Range("I8").NumberFormat = "$#,##0_);[Red]($#,##0)"

# Rows that carried a half-dollar amount -> "#,##0.00"
$ws.Range("I33").Value = 1187.5
$ws.Range("I34").Value = 1187.5
$ws.Range("I33:I34").NumberFormat = "#,##0.00"

# Remaining rows that become plain numbers with the default General format
$ws.Range("I69").Value = 120
$ws.Range("I79").Value = 525
$ws.Range("I88").Value = 525

# Leave the selection where the editor ended up (near the top of the newly
# reformatted Earnings column) instead of the original far-down J3 selection.
[void]$ws.Range("I10").Select()
